# Remove four obsolete / duplicated complaint rows from the "Optical_Power"
# sheet. Deleting them shifts every following row up, which reproduces the
# renumbering seen across rows 24-48 in the target diff (dimension becomes
# A1:L48 instead of A1:L52).
#
# Rows removed (original row numbers, identified by their unique "Caso"/OT
# values before any shifting occurs):
#   24 -> Caso 2225       / QUINQUELA MARTIN, BENITO 1282
#   26 -> Caso 4870       / ARAOZ DE LAMADRID, GREGORIO, GRAL. 283
#   35 -> Caso 5778       / GUEMES 3772
#   39 -> Caso 807168088  / Brandsen 1700
#
# They are deleted from the bottom up so that earlier row indices stay valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Delete()
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(24).Delete()
